$d = $word.ActiveDocument

# Locate the title text "Project-plan-v0.1"
$rng = $d.Content
$found = $rng.Find.Execute("Project-plan-v0.1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Replace the trailing "1" with "2"
    $lastCharRng = $d.Range($rng.End - 1, $rng.End)
    $lastCharRng.Text = "2"

    # Force the new "2" character into its own run (matching the author's
    # edit, which produced two adjacent runs with identical formatting)
    # by toggling a character property off and back on.
    $newCharRng = $d.Range($rng.End - 1, $rng.End)
    $newCharRng.Font.Bold = $false
    $newCharRng.Font.Bold = $true
}
